# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted as row 32 (pushing the
# existing rows 32-75 down to 33-76, i.e. a classic "insert row + shift
# down" edit). The sheet's used range grows from A1:R75 to A1:R76.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 32. This shifts the existing
# rows 32:75 down to 33:76, carrying their values/styles with them
# (Excel's native Insert-row behaviour), so no other row needs to be
# touched manually.
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new observation.
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 45079
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112026
$ws.Range("G32").Value = "Haba"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 30
$ws.Range("K32").Value = 15000
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = 15000
$ws.Range("N32").Value = '$/saco 25 kilos'
$ws.Range("O32").Value = "Provincia de Diguillín"
$ws.Range("P32").Value = 600
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
